$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Costos Legos"

$ws.Range("A1").Value = "ITEM"
$ws.Range("A2").Value = "Bases 12x12 para diseño"
$ws.Range("A3").Value = "Bloques Rojos"
$ws.Range("A4").Value = "Bloques Verdes"
$ws.Range("A5").Value = "Bloques Amarillos"
$ws.Range("A6").Value = "Bloques Azules"
$ws.Range("A7").Value = "TOTAL "
$ws.Range("A11").Value = "# Grupos"
$ws.Columns.Item(1).AutoFit() | Out-Null
Write-Host "done"
